$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing row 2 record down to row 3 (with updated Name/Email),
# then write the two new records into row 2 and row 4 -- without using
# Rows.Insert (which would drag row 1's header formatting down).

# Row 3: former row 2 record, Name/Email updated
$ws.Cells.Item(3, 1).Value = 631886740
$ws.Cells.Item(3, 2).Value = "Ikki maru"
$ws.Cells.Item(3, 3).Value = 992907510905
$ws.Cells.Item(3, 4).Value = "ikki@maru.com"
$ws.Cells.Item(3, 5).Value = "Ismat ."
$ws.Cells.Item(3, 6).Value = "EN"

# Row 2: new record
$ws.Cells.Item(2, 1).Value = 5547528084
$ws.Cells.Item(2, 2).Value = "Хушдил Саидов"
$ws.Cells.Item(2, 3).Value = 79177131361
$ws.Cells.Item(2, 4).Value = "xuwdi@mail.ru"
$ws.Cells.Item(2, 5).Value = "xuwdil None"
$ws.Cells.Item(2, 6).Value = "RU"

# Row 4: new record
$ws.Cells.Item(4, 1).Value = 974794263
$ws.Cells.Item(4, 2).Value = "Гульдартабакова Гульдартабакнукрахуросонабегим"
$ws.Cells.Item(4, 3).Value = 992938636344
$ws.Cells.Item(4, 4).Value = "thesarboz@gmail.com"
$ws.Cells.Item(4, 5).Value = "Buzurgmehr Abdulloev"
$ws.Cells.Item(4, 6).Value = "RU"
